# Update the Metadata sheet: Date and Count values
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-01-05T10:12:51-05:00"
$meta.Range("B21").Value = "'20"

# Update the Concepts sheet with the full disease list
$concepts = $wb.Worksheets.Item("Concepts")

$diseases = @(
    @("AN", "Anthrax"),
    @("WN", "West Nile"),
    @("BO", "Botulism"),
    @("CH", "Cholera"),
    @("CO", "Coronavirus"),
    @("DE", "Denque"),
    @("EB", "Ebola"),
    @("GO", "Gonorrhea"),
    @("ME", "Measles"),
    @("MA", "Malaria"),
    @("MB", "Marburg"),
    @("MN", "Meningococcal disease"),
    @("MP", "Mpox"),
    @("RA", "Rabies"),
    @("RU", "Rubella"),
    @("SP", "Smallpox"),
    @("SY", "Syphillis"),
    @("TB", "Tuberculosis"),
    @("YE", "Yellow Fever"),
    @("ZI", "Zika")
)

# Extend formatting (border/fill/alignment) of the existing data rows down
# through the new rows before writing values, so the new rows match the
# look of rows 2-3 (style index 2).
$concepts.Range("A2:D2").Copy()
$concepts.Range("A4:D21").PasteSpecial(-4122)

$row = 2
foreach ($d in $diseases) {
    $concepts.Cells.Item($row, 1).Value = "'1"
    $concepts.Cells.Item($row, 2).Value = $d[0]
    $concepts.Cells.Item($row, 3).Value = $d[1]
    $row = $row + 1
}
